# Fruta / hortaliza, semanal
# Rotate the weekly price records in rows 2-4: the record that was in row 2
# moves to row 3, row 3's record moves to row 4, and row 4's record moves
# up to row 2 (cyclic shift), updating columns D, L, M, N, O, P, R, S.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot current ("before") values for the columns that rotate.
$cols = @("D", "L", "M", "N", "O", "P", "R", "S")

$row2 = @{}
$row3 = @{}
$row4 = @{}
foreach ($col in $cols) {
    $row2[$col] = $ws.Range("${col}2").Value2
    $row3[$col] = $ws.Range("${col}3").Value2
    $row4[$col] = $ws.Range("${col}4").Value2
}

# Apply the rotation: new row2 = old row4, new row3 = old row2, new row4 = old row3
foreach ($col in $cols) {
    $ws.Range("${col}2").Value2 = $row4[$col]
    $ws.Range("${col}3").Value2 = $row2[$col]
    $ws.Range("${col}4").Value2 = $row3[$col]
}
